# Adds two new data rows (9 and 10) to Sheet1, extending the log_results
# table with two more "linear regression" / "1 row lookback" runs.
# This mirrors the "added support for dataframe to png" commit, which
# appended additional rows of simulation results to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "linear regression"
$ws.Cells.Item(9, 3).Value = "1 row lookback"
$ws.Cells.Item(9, 4).Value = 0.00000006538559915725273
$ws.Cells.Item(9, 5).Value = 0.0002031198964687064
$ws.Cells.Item(9, 6).Value = 103
$ws.Cells.Item(9, 7).Value = -0.000001947856617334764
$ws.Cells.Item(9, 8).Value = 0.0002569490170571953
$ws.Cells.Item(9, 9).Value = -0.0007576942443847656
$ws.Cells.Item(9, 10).Value = -0.0001809597015380859
$ws.Cells.Item(9, 11).Value = -0.000001311302185058594
$ws.Cells.Item(9, 12).Value = 0.0001505613327026367
$ws.Cells.Item(9, 13).Value = 0.0005565881729125977
$ws.Cells.Item(9, 14).Value = 0.0000004453056874353933
$ws.Cells.Item(9, 15).Value = 0.00000006538559915725273
$ws.Cells.Item(9, 16).Value = 0.0004266728064976633
$ws.Cells.Item(9, 17).Value = 0.0002031198964687064
$ws.Cells.Item(9, 18).Value = 0.00003986865858195155
$ws.Cells.Item(9, 19).Value = 0.5339805825242718
$ws.Cells.Item(9, 20).Value = 0
$ws.Cells.Item(9, 21).Value = 0.9636363636363636
$ws.Cells.Item(9, 22).Value = 0
$ws.Cells.Item(9, 23).Value = 0
$ws.Cells.Item(9, 24).Value = 0
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "linear regression"
$ws.Cells.Item(10, 3).Value = "1 row lookback"
$ws.Cells.Item(10, 4).Value = 0.00000006549666409227939
$ws.Cells.Item(10, 5).Value = 0.0002028722228715196
$ws.Cells.Item(10, 6).Value = 103
$ws.Cells.Item(10, 7).Value = -0.00002098430741170887
$ws.Cells.Item(10, 8).Value = 0.0002563086454756558
$ws.Cells.Item(10, 9).Value = -0.0007729530334472656
$ws.Cells.Item(10, 10).Value = -0.0001891851425170898
$ws.Cells.Item(10, 11).Value = -0.00002110004425048828
$ws.Cells.Item(10, 12).Value = 0.0001276731491088867
$ws.Cells.Item(10, 13).Value = 0.0005413293838500977
$ws.Cells.Item(10, 14).Value = 0.0000004451385393622331
$ws.Cells.Item(10, 15).Value = 0.00000006549666409227939
$ws.Cells.Item(10, 16).Value = 0.0004268484772183001
$ws.Cells.Item(10, 17).Value = 0.0002028722228715196
$ws.Cells.Item(10, 18).Value = 0.00003986865858195155
$ws.Cells.Item(10, 19).Value = 0.5242718446601942
$ws.Cells.Item(10, 20).Value = 0
$ws.Cells.Item(10, 21).Value = 0.9814814814814815
$ws.Cells.Item(10, 22).Value = 0
$ws.Cells.Item(10, 23).Value = 0
$ws.Cells.Item(10, 24).Value = 0

# Copy the index-column format (bold, centered, bordered) used on A2:A8
# down onto the two new index cells, same as every other row in column A.
$ws.Range("A8").Copy()
$ws.Range("A9:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
